$wb = $excel.ActiveWorkbook
$wsIB = $wb.Worksheets.Item("Initial Buys")
$ws15 = $wb.Worksheets.Item("2015")
$ws16 = $wb.Worksheets.Item("2016")
$ws17 = $wb.Worksheets.Item("2017")

# -----------------------------------------------------------------
# 1. "2017" sheet: insert a new row 11 (shifts old totals row to 12)
#    and populate it with the data that used to live in row 10
#    (the UL / Unilever holding), since row 10 becomes the new TGT
#    holding.
# -----------------------------------------------------------------
$ws17.Rows("11:11").Insert()

# Copy the formatting from row 9 (which uses the same banding as the
# old row 10 did before the insert) onto the freshly inserted row 11.
$ws17.Range("A9:V9").Copy()
$ws17.Range("A11:V11").PasteSpecial(-4122)
$ws17.Application.CutCopyMode = $false

$ws17.Range("A11").Value = "UL"
$ws17.Range("B11").Value = "Unilever"
$ws17.Range("C11").Value = "Consumer Staples"
$ws17.Range("D11").Value = 15.102
$ws17.Range("E11").Formula = "=V11/D11"
$ws17.Range("F11").Formula = "='Initial Buys'!Z40"
$ws17.Range("G11").Value = 0
$ws17.Range("H11").Value = 0
$ws17.Range("I11").Value = 0
$ws17.Range("J11").Value = 5.08
$ws17.Range("K11").Value = 0
$ws17.Range("L11").Value = 0
$ws17.Range("M11").Value = 0
$ws17.Range("N11").Value = 0
$ws17.Range("O11").Value = 0
$ws17.Range("P11").Value = 0
$ws17.Range("Q11").Value = 0
$ws17.Range("R11").Value = 0
$ws17.Range("S11").Value = 0
$ws17.Range("T11").Formula = "=SUM(H11:S11)"
$ws17.Range("U11").Formula = "=SUM(G11:S11)"
$ws17.Range("V11").Formula = "=SUM(F11,U11)"

# -----------------------------------------------------------------
# 2. "Initial Buys" sheet: add the new TGT buy block in columns AB/AC
# -----------------------------------------------------------------
$wsIB.Range("AB1").Value = "Date"
$wsIB.Range("Z1").Copy()
$wsIB.Range("AC1").PasteSpecial(-4122)
$wsIB.Application.CutCopyMode = $false
$wsIB.Range("AC1").Value = "TGT"

$wsIB.Range("AB2").Value = 42818
$wsIB.Range("AC2").Value = 797.25
$wsIB.Range("AC40").Formula = "=SUM(AC2:AC39)"

# -----------------------------------------------------------------
# 3. "2017" sheet: row 10 now represents the new TGT holding
# -----------------------------------------------------------------
$ws17.Range("A10").Value = "TGT"
$ws17.Range("B10").Value = "Target Corp"
$ws17.Range("C10").Value = "Consumer Discretionary"
$ws17.Range("D10").Value = 15
$ws17.Range("F10").Formula = "='Initial Buys'!AC40"
$ws17.Range("J10").Value = 0

# -----------------------------------------------------------------
# 4. "2017" sheet: rename the AT&T holding (row 9) and tweak its
#    industry label
# -----------------------------------------------------------------
$ws17.Range("B9").Value = "AT&T Inc"
$ws17.Range("C9").Value = "Telemunication"

# -----------------------------------------------------------------
# 5. "2017" sheet: updated share counts / dividend entries for a few
#    existing holdings (EMR row3, ADM row5, Shell row8)
# -----------------------------------------------------------------
$ws17.Range("D3").Value = 14.488
$ws17.Range("D5").Value = 40.382
$ws17.Range("K5").Value = 14.81
$ws17.Range("D8").Value = 19.052
$ws17.Range("J8").Value = 17.59

# -----------------------------------------------------------------
# 6. Column width tweak on "2017" sheet (Company Name column)
# -----------------------------------------------------------------
$ws17.Columns("C").ColumnWidth = 22.5703125

# -----------------------------------------------------------------
# 7. Sheet view / selection changes
# -----------------------------------------------------------------
$wsIB.Application.Goto($wsIB.Range("A22"), $true)
$wsIB.Range("Q33").Select()

$ws15.Application.Goto($ws15.Range("A1"), $true)
$ws15.Range("D15").Select()

$ws17.Range("K22").Select()
